# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" message on sheet "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rngA1 = $wsHoja1.Range("A1")
$oldText = $rngA1.Value()
$newText = $oldText.Replace("1000 Bs = 4.61 = 18295.3 pesos", "1000 Bs = 4.53 = 17994.97 pesos")
$newText = $newText.Replace("18295.3 pesos = 4.6 = 942.23 Bs", "17994.97 pesos = 4.51 = 951.15 Bs")
$rngA1.Value = $newText

# --- Update the rate cells on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 220.79
$wsTasas.Range("O10").Value = 3973.11
$wsTasas.Range("N12").Value = 3992
$wsTasas.Range("O12").Value = 211.003
